$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint 1 Example")

# Enter the new day-3 input values (columns E) that drive the recalculation
# of the burndown/cumulative-flow formulas and chart caches.
$ws.Range("E13").Value = 85
$ws.Range("E17").Value = 22
$ws.Range("E21").Value = 5
$ws.Range("E24").Value = 10

# Update the active view/selection on the sheet to match the saved state.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("E23").Select()
